$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header cell style (bold, bordered, centered) from H1
# onto the new header cells I1 and J1 so they match the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for rows 2 and 3
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8
